# regression2: resolved CI calc problems and some ANOVA values
# Adds b0, b1, b2, b11, b22, b12 design-matrix columns (F:K) derived from
# the existing coded factor columns A (x1) and C (x2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: F1:K1
$ws.Range("F1").Value = "b0"
$ws.Range("G1").Value = "b1"
$ws.Range("H1").Value = "b2"
$ws.Range("I1").Value = "b11"
$ws.Range("J1").Value = "b22"
$ws.Range("K1").Value = "b12"

# Data rows 2..10
For ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
    $ws.Cells.Item($r, 7).Formula = "=A" + $r
    $ws.Cells.Item($r, 8).Formula = "=C" + $r
    $ws.Cells.Item($r, 9).Formula = "=G" + $r + "^2"
    $ws.Cells.Item($r, 10).Formula = "=H" + $r + "^2"
    $ws.Cells.Item($r, 11).Formula = "=G" + $r + "*H" + $r
}

$ws.Range("F11").Select() | Out-Null
